# Add a "team record" (Wins / Losses / Ties) block to the right of the
# existing roster columns (data runs through column AC on rows 1-66).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels in row 1, columns AD:AF.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the other header cells (A1:AC1): bold font, thin box
# border, centered horizontally and top-aligned vertically.
$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous

# Every player row (2-66) shares the team's overall 2021 record.
$ws.Range("AD2:AD66").Value = 77
$ws.Range("AE2:AE66").Value = 85
$ws.Range("AF2:AF66").Value = 0
